$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12472.667
$ws.Range("I62").Value = 1531.75
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 1531.75
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -907.75
$ws.Range("N62").Value = -101248

$ws.Range("H65").Value = 12472.667
$ws.Range("I65").Value = 1531.75
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 7658.75
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -4538.75
$ws.Range("N65").Value = -506240

$ws.Range("H107").Value = 837.75
$ws.Range("I107").Value = 1057.85
$ws.Range("J107").Value = 287.5
$ws.Range("K107").Value = 1057.85
$ws.Range("L107").Value = 287.5
$ws.Range("M107").Value = 862.1500000000001
$ws.Range("N107").Value = -4127.5

$ws.Range("H111").Value = 3223.1
$ws.Range("I111").Value = 1885.5714
$ws.Range("J111").Value = 6344
$ws.Range("K111").Value = 5656.7142
$ws.Range("L111").Value = 19032
$ws.Range("M111").Value = -2589.7142
$ws.Range("N111").Value = -25166

$ws.Range("H129").Value = 813.63336
$ws.Range("I129").Value = 634.9286
$ws.Range("J129").Value = 970
$ws.Range("K129").Value = 1904.7858
$ws.Range("L129").Value = 2910
$ws.Range("M129").Value = 3095.2142
$ws.Range("N129").Value = -12910

$ws.Range("H137").Value = 2316.5557
$ws.Range("I137").Value = 1630.3889
$ws.Range("J137").Value = 3688.889
$ws.Range("K137").Value = 4891.1667
$ws.Range("L137").Value = 11066.667
$ws.Range("M137").Value = -2341.1667
$ws.Range("N137").Value = -16166.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3598.56
$ws.Range("I45").Value = 2761.6428
$ws.Range("J45").Value = 4663.727
$ws.Range("K45").Value = 2761.6428
$ws.Range("L45").Value = 4663.727
$ws.Range("M45").Value = -2384.6428
$ws.Range("N45").Value = -5417.727

$ws.Range("H61").Value = 2236.95
$ws.Range("I61").Value = 2053.6875
$ws.Range("J61").Value = 2970
$ws.Range("K61").Value = 2053.6875
$ws.Range("L61").Value = 2970
$ws.Range("M61").Value = -1841.6875
$ws.Range("N61").Value = -3394

$ws.Range("H63").Value = 3222.111
$ws.Range("I63").Value = 1785.4286
$ws.Range("K63").Value = 1785.4286
$ws.Range("M63").Value = -1099.4286

$ws.Range("H66").Value = 3222.111
$ws.Range("I66").Value = 1785.4286
$ws.Range("K66").Value = 8927.143
$ws.Range("M66").Value = -5495.143

$ws.Range("H74").Value = 1242
$ws.Range("I74").Value = 1251.2
$ws.Range("K74").Value = 1251.2
$ws.Range("M74").Value = -377.2

$ws.Range("H77").Value = 1242
$ws.Range("I77").Value = 1251.2
$ws.Range("K77").Value = 6256
$ws.Range("M77").Value = -1888

$ws.Range("H122").Value = 1816.909
$ws.Range("I122").Value = 1754.5555
$ws.Range("J122").Value = 2097.5
$ws.Range("K122").Value = 5263.666499999999
$ws.Range("L122").Value = 6292.5
$ws.Range("M122").Value = -2813.666499999999
$ws.Range("N122").Value = -11192.5

$ws.Range("H136").Value = 2236.95
$ws.Range("I136").Value = 2053.6875
$ws.Range("J136").Value = 2970
$ws.Range("K136").Value = 6161.0625
$ws.Range("L136").Value = 8910
$ws.Range("M136").Value = -3611.0625
$ws.Range("N136").Value = -14010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24790.422
$ws.Range("J82").Value = 32845.08
$ws.Range("L82").Value = 32845.08
$ws.Range("N82").Value = -33611.08

$ws.Range("H85").Value = 24790.422
$ws.Range("J85").Value = 32845.08
$ws.Range("L85").Value = 32845.08
$ws.Range("N85").Value = -35497.08

$ws.Range("H134").Value = 1335.7826
$ws.Range("I134").Value = 1282.0476
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 3846.142800000001
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = -1311.142800000001
$ws.Range("N134").Value = -10770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2585.8823
$ws.Range("I31").Value = 1828.1333
$ws.Range("J31").Value = 3184.1052
$ws.Range("K31").Value = 1828.1333
$ws.Range("L31").Value = 3184.1052
$ws.Range("M31").Value = -1533.1333
$ws.Range("N31").Value = -3774.1052

$ws.Range("H34").Value = 2585.8823
$ws.Range("I34").Value = 1828.1333
$ws.Range("J34").Value = 3184.1052
$ws.Range("K34").Value = 1828.1333
$ws.Range("L34").Value = 3184.1052
$ws.Range("M34").Value = -1626.1333
$ws.Range("N34").Value = -3588.1052

$ws.Range("H50").Value = 8945.143
$ws.Range("J50").Value = 8945.143
$ws.Range("L50").Value = 8945.143
$ws.Range("N50").Value = -10195.143

$ws.Range("H60").Value = 22514
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H122").Value = 2076.4375
$ws.Range("I122").Value = 1283.909
$ws.Range("J122").Value = 3820
$ws.Range("K122").Value = 3851.727
$ws.Range("L122").Value = 11460
$ws.Range("M122").Value = -1401.727
$ws.Range("N122").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28427.273
$ws.Range("I70").Value = 33611.11
$ws.Range("J70").Value = 5100
$ws.Range("K70").Value = 33611.11
$ws.Range("L70").Value = 5100
$ws.Range("M70").Value = -33341.11
$ws.Range("N70").Value = -5640

$ws.Range("H73").Value = 28427.273
$ws.Range("I73").Value = 33611.11
$ws.Range("J73").Value = 5100
$ws.Range("K73").Value = 33611.11
$ws.Range("L73").Value = 5100
$ws.Range("M73").Value = -32675.11
$ws.Range("N73").Value = -6972

$ws.Range("H122").Value = 1685
$ws.Range("I122").Value = 1701.8889
$ws.Range("K122").Value = 5105.6667
$ws.Range("M122").Value = -2655.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1977.1111
$ws.Range("I40").Value = 1827.7142
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1827.7142
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1691.7142
$ws.Range("N40").Value = -2772

$ws.Range("H55").Value = 463.1111
$ws.Range("I55").Value = 137.66667
$ws.Range("J55").Value = 723.4666999999999
$ws.Range("K55").Value = 137.66667
$ws.Range("L55").Value = 723.4666999999999
$ws.Range("M55").Value = 35.33332999999999
$ws.Range("N55").Value = -1069.4667

$ws.Range("H136").Value = 2306.8572
$ws.Range("I136").Value = 1880.2222
$ws.Range("J136").Value = 4866.6665
$ws.Range("K136").Value = 5640.6666
$ws.Range("L136").Value = 14599.9995
$ws.Range("M136").Value = -3090.6666
$ws.Range("N136").Value = -19699.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1209.2162
$ws.Range("I107").Value = 861
$ws.Range("J107").Value = 1666.25
$ws.Range("K107").Value = 2583
$ws.Range("L107").Value = 4998.75
$ws.Range("M107").Value = -663
$ws.Range("N107").Value = -8838.75

$ws.Range("H109").Value = 28788.5
$ws.Range("J109").Value = 28788.5
$ws.Range("L109").Value = 28788.5
$ws.Range("N109").Value = -31562.5
